$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.435.76"
$ws.Range("E2").Value = "  -2.11%  "
$ws.Range("D3").Value = "2.276.38"
$ws.Range("E3").Value = "  -2.23%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "298.91"
$ws.Range("E5").Value = "  -1.74%  "
$ws.Range("D6").Value = "95.29"
$ws.Range("E6").Value = "  -5.55%  "
$ws.Range("E7").Value = "  -1.90%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -3.37%  "
$ws.Range("D10").Value = "33.06"
$ws.Range("E10").Value = "  -6.04%  "
$ws.Range("D11").Value = "0.0788"
$ws.Range("E11").Value = "  -1.09%  "
$ws.Range("D12").Value = "48.96"
$ws.Range("E12").Value = "  -3.94%  "
$ws.Range("E13").Value = "  +1.91%  "
$ws.Range("D14").Value = "16.61"
$ws.Range("E14").Value = "  +6.40%  "
$ws.Range("E15").Value = "  -1.01%  "
$ws.Range("D16").Value = "2.630.66"
$ws.Range("E16").Value = "  -2.27%  "
$ws.Range("D17").Value = "2.286.47"
$ws.Range("E17").Value = "  -0.90%  "
$ws.Range("E18").Value = "  -1.57%  "
$ws.Range("D19").Value = "42.342.43"
$ws.Range("E19").Value = "  -2.15%  "
$ws.Range("E20").Value = "  -1.88%  "
$ws.Range("D21").Value = "11.37"
$ws.Range("E21").Value = "  -3.79%  "
$ws.Range("E22").Value = "  -2.51%  "
$ws.Range("D23").Value = "66.65"
$ws.Range("E23").Value = "  -2.15%  "
$ws.Range("D24").Value = "235.54"
$ws.Range("E24").Value = "  -0.74%  "
$ws.Range("E25").Value = "  -1.51%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  -3.51%  "
$ws.Range("D28").Value = "24.15"
$ws.Range("E28").Value = "  -3.06%  "
$ws.Range("E29").Value = "  -0.83%  "
$ws.Range("D30").Value = "166.32"
$ws.Range("E30").Value = "  +0.55%  "
$ws.Range("D31").Value = "33.41"
$ws.Range("D32").Value = "9.05"
$ws.Range("E32").Value = "  -1.80%  "
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("E34").Value = "  +3.35%  "
$ws.Range("D35").Value = "4.91"
$ws.Range("E35").Value = "  -2.76%  "
$ws.Range("D36").Value = "2.41"
$ws.Range("E36").Value = "  -0.32%  "
$ws.Range("D37").Value = "16.70"
$ws.Range("E37").Value = "  -0.69%  "
$ws.Range("E38").Value = "  -2.40%  "
$ws.Range("E39").Value = "  -4.13%  "
$ws.Range("D40").Value = "0.0990"
$ws.Range("E40").Value = "  -3.53%  "
$ws.Range("E41").Value = "  -1.85%  "
$ws.Range("E42").Value = "  -5.16%  "
$ws.Range("D43").Value = "2.36"
$ws.Range("E43").Value = "  -2.86%  "
$ws.Range("D44").Value = "1.950.25"
$ws.Range("E44").Value = "  -1.15%  "
$ws.Range("D45").Value = "0.0277"
$ws.Range("E45").Value = "  -2.13%  "
$ws.Range("E46").Value = "  -3.47%  "
$ws.Range("D47").Value = "17.26"
$ws.Range("E47").Value = "  -6.91%  "
$ws.Range("E48").Value = "  -4.37%  "
$ws.Range("D49").Value = "2.503.47"
$ws.Range("E49").Value = "  -1.94%  "
$ws.Range("D50").Value = "52.18"
$ws.Range("E50").Value = "  -6.52%  "
$ws.Range("D51").Value = "2.73"
$ws.Range("E51").Value = "  -2.27%  "
